# Generate Report for Archive
#
# The localization status report is regenerated: the rows describing the
# files "071d1400-...", "e3558659-..." and "fa9b0075-..." are re-ordered
# (fa9b0075 moves up to the first of the three slots, and its status flips
# from "Ready for handoff" to "In Translation"; 071d1400 and e3558659 shift
# down one slot each, keeping their own data). This happens identically on
# the "Overview" sheet (rows 5-7, columns A-D) and on each per-language
# sheet "zh-cn" / "de-de" (rows 5-7, columns A-E).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: columns A (file), B/C (status x2), D (handoff date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A5").Value = "fa9b0075-b570-4301-a336-d2c339cecb02.md"
$wsOverview.Range("B5").Value = "In Translation"
$wsOverview.Range("C5").Value = "In Translation"
$wsOverview.Range("D5").Value = "2016-41-13 06:41:32"

$wsOverview.Range("A6").Value = "071d1400-152d-4846-8ff8-8ab201631f57.md"
$wsOverview.Range("B6").Value = "Ready for handoff"
$wsOverview.Range("C6").Value = "Ready for handoff"
$wsOverview.Range("D6").Value = "2016-39-13 06:39:45"

$wsOverview.Range("A7").Value = "e3558659-c60c-420c-9c1a-ef4fc13ba77e.md"
$wsOverview.Range("B7").Value = "Ready for handoff"
$wsOverview.Range("C7").Value = "Ready for handoff"
$wsOverview.Range("D7").Value = "2016-41-13 06:41:49"

# ---------------------------------------------------------------------
# zh-cn sheet: columns A (file), B (ext), C (status), D (target file), E (date)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A5").Value = "fa9b0075-b570-4301-a336-d2c339cecb02.md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("D5").Value = "fa9b0075-b570-4301-a336-d2c339cecb02.eac8b026d1bd6452827445c0d29d225c9d4b43ba.zh-cn.xlf"
$wsZhCn.Range("E5").Value = "2016-03-13 06:41:29"

$wsZhCn.Range("A6").Value = "071d1400-152d-4846-8ff8-8ab201631f57.md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("D6").Value = "071d1400-152d-4846-8ff8-8ab201631f57.a69c322b52248a332fe2d0ea4529f83daa92a0d8.zh-cn.xlf"
$wsZhCn.Range("E6").Value = "2016-03-13 06:39:41"

$wsZhCn.Range("A7").Value = "e3558659-c60c-420c-9c1a-ef4fc13ba77e.md"
$wsZhCn.Range("C7").Value = "Ready for handoff"
$wsZhCn.Range("D7").Value = "e3558659-c60c-420c-9c1a-ef4fc13ba77e.4760866b4cb1116d9eb067328510f2ea2a09d5ca.zh-cn.xlf"
$wsZhCn.Range("E7").Value = "2016-03-13 06:41:46"

# ---------------------------------------------------------------------
# de-de sheet: columns A (file), B (ext), C (status), D (target file), E (date)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A5").Value = "fa9b0075-b570-4301-a336-d2c339cecb02.md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("D5").Value = "fa9b0075-b570-4301-a336-d2c339cecb02.eac8b026d1bd6452827445c0d29d225c9d4b43ba.de-de.xlf"
$wsDeDe.Range("E5").Value = "2016-03-13 06:41:32"

$wsDeDe.Range("A6").Value = "071d1400-152d-4846-8ff8-8ab201631f57.md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("D6").Value = "071d1400-152d-4846-8ff8-8ab201631f57.a69c322b52248a332fe2d0ea4529f83daa92a0d8.de-de.xlf"
$wsDeDe.Range("E6").Value = "2016-03-13 06:39:45"

$wsDeDe.Range("A7").Value = "e3558659-c60c-420c-9c1a-ef4fc13ba77e.md"
$wsDeDe.Range("C7").Value = "Ready for handoff"
$wsDeDe.Range("D7").Value = "e3558659-c60c-420c-9c1a-ef4fc13ba77e.4760866b4cb1116d9eb067328510f2ea2a09d5ca.de-de.xlf"
$wsDeDe.Range("E7").Value = "2016-03-13 06:41:49"
